$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns stay text-typed,
# matching the source workbook where every D/E data cell is inline string text
# (Price values like "1.00" / "3.340.34" are NOT numbers in this sheet).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "65.061.43"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "3.348.23"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "558.38"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").Value = "173.62"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "3.337.62"
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").Value = "0.629"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "52.86"
$ws.Range("E12").Value = "  -3.70%  "
$ws.Range("D13").Value = "0.0000275"
$ws.Range("D14").Value = "9.13"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "3.888.82"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "0.120"
$ws.Range("E16").Value = "  +1.04%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "18.10"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.325.80"
$ws.Range("E18").Value = "  -2.58%  "
$ws.Range("D19").Value = "64.992.03"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").Value = "11.72"
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").Value = "0.990"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").Value = "480.96"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").Value = "4.92"
$ws.Range("E23").Value = "  -4.11%  "
$ws.Range("D24").Value = "89.65"
$ws.Range("E24").Value = "  +3.53%  "
$ws.Range("D25").Value = "14.26"
$ws.Range("E25").Value = "  +5.23%  "
$ws.Range("D26").Value = "4.07"
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("D27").Value = "2.88"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "10.51"
$ws.Range("E28").Value = "  -3.52%  "
$ws.Range("D29").Value = "8.64"
$ws.Range("E29").Value = "  -2.72%  "
$ws.Range("D30").Value = "31.07"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "6.50"
$ws.Range("E31").Value = "  -2.70%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "11.37"
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "571.89"
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "61.58"
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "3.63"
$ws.Range("E37").Value = "  +3.96%  "
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").Value = "35.43"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").Value = "0.372"
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("D41").Value = "0.0₃0732"
$ws.Range("E41").Value = "  -2.77%  "
$ws.Range("D42").Value = "3.097.13"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").Value = "2.78"
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("D44").Value = "0.0413"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "3.17"
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.133"
$ws.Range("E46").Value = "  -1.28%  "
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  -2.75%  "
$ws.Range("D48").Value = "0.997"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").Value = "140.98"
$ws.Range("E49").Value = "  +3.41%  "
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("D51").Value = "8.37"
$ws.Range("E51").Value = "  +0.33%  "
